# Fix generate evaluation form
#
# Two of the document's three tables are missing <w:tblLayout w:type="fixed"/>
# in their <w:tblPr> (it sits right after </w:tblBorders> and before
# <w:tblLook .../>). The third table already has it. Setting
# Table.AllowAutoFit = $false is the Word object-model equivalent of
# writing <w:tblLayout w:type="fixed"/>, so apply it to the two tables
# that currently lack it (leave the one that already has it alone).

$d = $word.ActiveDocument

for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    if ($t.AllowAutoFit -ne $false) {
        $t.AllowAutoFit = $false
    }
}
